# Reorder the "Recorded By" (column G) names in each data row.
# For every row whose "Recorded By" cell holds more than one comma-separated
# name/address, the last entry is moved to the front (a right-rotation by
# one position). Cells holding a single value are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow   = $usedRange.Rows.Count
$lastCol   = $usedRange.Columns.Count

# Locate the "Recorded By" column dynamically (falls back to column 7).
$recordedByCol = 0
for ($c = 1; $c -le $lastCol; $c++) {
    if ($ws.Cells.Item(1, $c).Text -eq "Recorded By") {
        $recordedByCol = $c
        break
    }
}
if ($recordedByCol -eq 0) {
    $recordedByCol = 7
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $recordedByCol)
    $val = $cell.Text

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Length -gt 1) {
            $newParts = @($parts[-1]) + $parts[0..($parts.Length - 2)]
            $newVal = $newParts -join ", "
            $cell.Value = $newVal
        }
    }
}
